$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1984.1177
$ws.Range("I40").Value = 1613.0769
$ws.Range("J40").Value = 2213.8096
$ws.Range("K40").Value = 1613.0769
$ws.Range("L40").Value = 2213.8096
$ws.Range("M40").Value = -1438.0769
$ws.Range("N40").Value = -2563.8096
$ws.Range("H43").Value = 3165.2856
$ws.Range("I43").Value = 5766.6665
$ws.Range("J43").Value = 1214.25
$ws.Range("K43").Value = 5766.6665
$ws.Range("L43").Value = 1214.25
$ws.Range("M43").Value = -5697.6665
$ws.Range("N43").Value = -1352.25
$ws.Range("H96").Value = 1381.3334
$ws.Range("I96").Value = 1153.2222
$ws.Range("J96").Value = 1723.5
$ws.Range("K96").Value = 3459.6666
$ws.Range("L96").Value = 5170.5
$ws.Range("M96").Value = -2086.6666
$ws.Range("N96").Value = -7916.5
$ws.Range("H112").Value = 1658.1818
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1658.1818
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4974.5454
$ws.Range("N112").Value = -7190.5454
$ws.Range("M112").ClearContents()
$ws.Range("H116").Value = 5278
$ws.Range("I116").Value = 5254.4443
$ws.Range("J116").Value = 5313.3335
$ws.Range("K116").Value = 5254.4443
$ws.Range("L116").Value = 5313.3335
$ws.Range("M116").Value = -1812.4443
$ws.Range("N116").Value = -12197.3335
$ws.Range("H132").Value = 2231.5938
$ws.Range("I132").Value = 2716.2632
$ws.Range("J132").Value = 1523.2307
$ws.Range("K132").Value = 8148.7896
$ws.Range("L132").Value = 4569.6921
$ws.Range("M132").Value = -5618.7896
$ws.Range("N132").Value = -9629.6921
$ws.Range("H139").Value = 52000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 52000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 52000
$ws.Range("N139").Value = -62280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8592.433999999999
$ws.Range("I32").Value = 7404.8335
$ws.Range("J32").Value = 13342.833
$ws.Range("K32").Value = 7404.8335
$ws.Range("L32").Value = 13342.833
$ws.Range("M32").Value = -7117.8335
$ws.Range("N32").Value = -13916.833
$ws.Range("H102").Value = 4971.4287
$ws.Range("I102").Value = 4971.4287
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4971.4287
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3349.4287
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 1105.5
$ws.Range("I110").Value = 1105.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1105.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 939.5
$ws.Range("H132").Value = 3900.3333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3900.3333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11700.9999
$ws.Range("N132").Value = -16760.9999
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 19285
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 19285
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 19285
$ws.Range("N81").Value = -21407
$ws.Range("H84").Value = 19285
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 19285
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 57855
$ws.Range("N84").Value = -68463
$ws.Range("H107").Value = 13446.75
$ws.Range("I107").Value = 2084.5557
$ws.Range("J107").Value = 47533.332
$ws.Range("K107").Value = 2084.5557
$ws.Range("L107").Value = 47533.332
$ws.Range("M107").Value = -164.5556999999999
$ws.Range("N107").Value = -51373.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 21500
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 21500
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 21500
$ws.Range("N26").Value = -22074
$ws.Range("H31").Value = 2338.8645
$ws.Range("I31").Value = 2361.255
$ws.Range("J31").Value = 2196.125
$ws.Range("K31").Value = 2361.255
$ws.Range("L31").Value = 2196.125
$ws.Range("M31").Value = -2066.255
$ws.Range("N31").Value = -2786.125
$ws.Range("H34").Value = 2338.8645
$ws.Range("I34").Value = 2361.255
$ws.Range("J34").Value = 2196.125
$ws.Range("K34").Value = 2361.255
$ws.Range("L34").Value = 2196.125
$ws.Range("M34").Value = -2159.255
$ws.Range("N34").Value = -2600.125
$ws.Range("H132").Value = 1698.08
$ws.Range("I132").Value = 1184.9375
$ws.Range("J132").Value = 2610.3333
$ws.Range("K132").Value = 3554.8125
$ws.Range("L132").Value = 7830.999899999999
$ws.Range("M132").Value = -1024.8125
$ws.Range("N132").Value = -12890.9999
$ws.Range("H134").Value = 33334744
$ws.Range("I134").Value = 1446.5
$ws.Range("J134").Value = 166667940
$ws.Range("K134").Value = 4339.5
$ws.Range("L134").Value = 500003820
$ws.Range("M134").Value = -1804.5
$ws.Range("N134").Value = -500008890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 767
$ws.Range("I5").Value = 633.1667
$ws.Range("J5").Value = 824.3570999999999
$ws.Range("K5").Value = 1899.5001
$ws.Range("L5").Value = 2473.0713
$ws.Range("M5").Value = -1787.5001
$ws.Range("N5").Value = -2697.0713
$ws.Range("H68").Value = 1592.3334
$ws.Range("I68").Value = 850
$ws.Range("J68").Value = 1804.4286
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 5413.2858
$ws.Range("M68").Value = -1739
$ws.Range("N68").Value = -7035.2858
$ws.Range("H71").Value = 1592.3334
$ws.Range("I71").Value = 850
$ws.Range("J71").Value = 1804.4286
$ws.Range("K71").Value = 7650
$ws.Range("L71").Value = 16239.8574
$ws.Range("M71").Value = -3594
$ws.Range("N71").Value = -24351.8574
$ws.Range("H122").Value = 701.1111
$ws.Range("I122").Value = 466.66666
$ws.Range("J122").Value = 768.0952
$ws.Range("K122").Value = 4199.99994
$ws.Range("L122").Value = 6912.8568
$ws.Range("M122").Value = -1749.99994
$ws.Range("N122").Value = -11812.8568
$ws.Range("H123").Value = 4078.625
$ws.Range("I123").Value = 910.3333
$ws.Range("J123").Value = 5979.6
$ws.Range("K123").Value = 2730.9999
$ws.Range("L123").Value = 17938.8
$ws.Range("M123").Value = -280.9998999999998
$ws.Range("N123").Value = -22838.8
$ws.Range("H129").Value = 3291.889
$ws.Range("I129").Value = 1314
$ws.Range("J129").Value = 4455.353
$ws.Range("K129").Value = 3942
$ws.Range("L129").Value = 13366.059
$ws.Range("M129").Value = 1058
$ws.Range("N129").Value = -23366.059
$ws.Range("H135").Value = 767
$ws.Range("I135").Value = 633.1667
$ws.Range("J135").Value = 824.3570999999999
$ws.Range("K135").Value = 5698.5003
$ws.Range("L135").Value = 7419.2139
$ws.Range("M135").Value = -3163.5003
$ws.Range("N135").Value = -12489.2139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1743.5294
$ws.Range("I102").Value = 1255.6364
$ws.Range("J102").Value = 2638
$ws.Range("K102").Value = 1255.6364
$ws.Range("L102").Value = 2638
$ws.Range("M102").Value = 366.3635999999999
$ws.Range("N102").Value = -5882
$ws.Range("H122").Value = 940813.3
$ws.Range("I122").Value = 1316608.9
$ws.Range("J122").Value = 1324.25
$ws.Range("K122").Value = 3949826.7
$ws.Range("L122").Value = 3972.75
$ws.Range("M122").Value = -3947376.7
$ws.Range("N122").Value = -8872.75
$ws.Range("H132").Value = 3880.842
$ws.Range("I132").Value = 4503
$ws.Range("J132").Value = 3428.3635
$ws.Range("K132").Value = 13509
$ws.Range("L132").Value = 10285.0905
$ws.Range("M132").Value = -10979
$ws.Range("N132").Value = -15345.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2485
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 2980
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 2980
$ws.Range("M2").Value = -888
$ws.Range("N2").Value = -3204
$ws.Range("H13").Value = 4000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4280
$ws.Range("H22").Value = 1456.25
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1573.6111
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1573.6111
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -2163.6111
$ws.Range("H27").Value = 1456.25
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1573.6111
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1573.6111
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1787.6111
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 2779.9
$ws.Range("I132").Value = 1844.4445
$ws.Range("J132").Value = 3545.2727
$ws.Range("K132").Value = 5533.333500000001
$ws.Range("L132").Value = 10635.8181
$ws.Range("M132").Value = -3003.333500000001
$ws.Range("N132").Value = -15695.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3713
$ws.Range("I62").Value = 3665.7778
$ws.Range("J62").Value = 3798
$ws.Range("K62").Value = 3665.7778
$ws.Range("L62").Value = 3798
$ws.Range("M62").Value = -3041.7778
$ws.Range("N62").Value = -5046
$ws.Range("H65").Value = 3713
$ws.Range("I65").Value = 3665.7778
$ws.Range("J65").Value = 3798
$ws.Range("K65").Value = 18328.889
$ws.Range("L65").Value = 18990
$ws.Range("M65").Value = -15208.889
$ws.Range("N65").Value = -25230
$ws.Range("H100").Value = 2880
$ws.Range("I100").Value = 2880
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5760
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5219
$ws.Range("H136").Value = 519.4400000000001
$ws.Range("I136").Value = 499.41666
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 1498.24998
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = 1051.75002
$ws.Range("N136").Value = -8100

"done"